$d = $word.ActiveDocument

$replacements = @(
    @("81×37=", "63×15="),
    @("67×11=", "61×61="),
    @("53×97=", "60×75="),
    @("40×26=", "23×93="),
    @("70×79=", "62×88="),
    @("99×95=", "49×99="),
    @("11×97=", "26×12="),
    @("19×89=", "99×73="),
    @("98×51=", "63×33="),
    @("69×23=", "34×75="),
    @("72×63=", "34×91="),
    @("96×38=", "20×18="),
    @("77×92=", "69×54="),
    @("63×44=", "91×77="),
    @("88×26=", "37×38="),
    @("41×20=", "24×31="),
    @("18×88=", "25×81="),
    @("46×32=", "54×72="),
    @("93×79=", "34×76="),
    @("78×51=", "66×55="),
    @("96×47=", "41×71="),
    @("56×73=", "66×11="),
    @("44×42=", "35×35="),
    @("50×93=", "61×62="),
    @("87×36=", "36×13=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
